# Legacy GSC export refresh: the oldest day (2025-10-18) drops off the
# rolling window, every later day's row shifts up one position, and the
# two newest days (2025-10-19 / 2025-10-20) don't have "Not indexed" /
# "Indexed" figures yet, so those land as blank text cells instead of
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the 2025-10-18 row entirely; everything below shifts up and the
# old last row (89) falls off, shrinking the used range to A1:D88.
$ws.Rows.Item(2).Delete()

# The two most recent dates (now rows 2 and 3) don't have "Not indexed"
# / "Indexed" counts yet in this export -- store them as empty text
# cells (not simply cleared cells).
$ws.Range("B2:C3").Value = "'"
$ws.Range("B2:C3").ClearFormats()
